$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 23; $r -le 245; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
